# Update Name of Algo
# Apply updated imputed values to the result data worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.2006
$ws.Range("C3").Value = -11.49929999999999
$ws.Range("A4").Value = -21.32050000000002
$ws.Range("B4").Value = 4.674200000000003
$ws.Range("C4").Value = -11.07989999999999
$ws.Range("B5").Value = 5.455399999999996
$ws.Range("E5").Value = 13.31809999999999
$ws.Range("A6").Value = -21.65080000000002
$ws.Range("B6").Value = 5.273899999999994
$ws.Range("A7").Value = -21.28740000000002
$ws.Range("A8").Value = -21.45740000000003
$ws.Range("B8").Value = 4.934
$ws.Range("C9").Value = -11.4373
$ws.Range("C11").Value = -14.40440000000001
$ws.Range("C14").Value = -11.28649999999999
$ws.Range("A16").Value = -21.46800000000002
$ws.Range("B16").Value = 5.427299999999998
$ws.Range("C18").Value = -14.68580000000001
$ws.Range("A20").Value = -22.60080000000003
$ws.Range("E20").Value = 13.30989999999999
$ws.Range("A21").Value = -20.5837
$ws.Range("B22").Value = 5.173900000000003
$ws.Range("C25").Value = -11.49189999999999
